# Update Data by bot, scripted by HH
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2020-09-30 00:00:00"

$ws.Range("O2").Value = -85733936.78
$ws.Range("P2").Value = -74.5799260948
$ws.Range("Q2").Value = 647013066.28
$ws.Range("R2").Value = 562.8364738384
$ws.Range("S2").Value = 60968172.97
$ws.Range("T2").Value = 53.0361955255
$ws.Range("U2").Value = -278755.71
$ws.Range("V2").Value = -0.2424895092
$ws.Range("W2").Value = 5260390
$ws.Range("X2").Value = 4.5760116958
$ws.Range("Y2").Value = 10316485.57
$ws.Range("Z2").Value = 8.9743077277
$ws.Range("AA2").Value = -28943088.87
$ws.Range("AB2").Value = -25.177584396
$ws.Range("AC2").Value = -114955781.36
$ws.Range("AD2").Value = 88.920052769
